$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells in column D store text-formatted numeric strings (inlineStr in the
# source data), so force Text number format before assigning the value to avoid
# Excel auto-converting the string into a numeric value.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D27", "D40", "D41", "D42", "D43", "D44", "D48")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Column D (Price) updates ---
$ws.Range("D2").Value = '244.04'
$ws.Range("D3").Value = '25.07'
$ws.Range("D4").Value = '5.185'
$ws.Range("D5").Value = '0.05747'
$ws.Range("D6").Value = '6.512'
$ws.Range("D8").Value = '0.8101'
$ws.Range("D9").Value = '0.8383'
$ws.Range("D10").Value = '0.1338'
$ws.Range("D11").Value = '0.06960'
$ws.Range("D12").Value = '0.03127'
$ws.Range("D13").Value = '0.02828'
$ws.Range("D14").Value = '0.09369'
$ws.Range("D15").Value = '0.001518'
$ws.Range("D16").Value = '0.0005999'
$ws.Range("D17").Value = '0.006223'
$ws.Range("D18").Value = '3.501'
$ws.Range("D19").Value = '2.092'
$ws.Range("D20").Value = '0.3175'
$ws.Range("D22").Value = '3.741'
$ws.Range("D23").Value = '0.04655'
$ws.Range("D24").Value = '0.1329'
$ws.Range("D25").Value = '0.001236'
$ws.Range("D27").Value = '0.00008699'
$ws.Range("D40").Value = '0.03609'
$ws.Range("D41").Value = '0.006299'
$ws.Range("D42").Value = '0.1050'
$ws.Range("D43").Value = '0.003000'
$ws.Range("D44").Value = '0.007321'
$ws.Range("D48").Value = '0.002283'

# --- Coin / Link / Volume(1h) updates (rows shifted due to a new coin insertion) ---
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("E15").Value = '14BitForexTokenBF'
$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("E16").Value = '15OneONEWorstin24h'
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("E17").Value = '16TigerCashTCH'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("E18").Value = '17LEOLEO'
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("E19").Value = '18BTSETokenBTSE'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("E20").Value = '19BitpandaEcosystemTokenBEST'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("E43").Value = '42CEJICEJIBestin24h'
